# Generate Report for Handback
# Marks the 0b47b63f... and c617b04f... localization jobs as handed back
# (in sync with en-US) for both the zh-cn and de-de targets, and records
# the generated target / handback file names + handback timestamps.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$srcMdA  = "0b47b63f-f910-4fa8-8c79-d4ab304f577c.md"
$srcMdB  = "c617b04f-6424-4b11-adb3-db1762f1a537.md"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfe93e339115c5a4322378e793c43696e8085dc6/e2e/0b47b63f-f910-4fa8-8c79-d4ab304f577c.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dfe93e339115c5a4322378e793c43696e8085dc6/e2e/c617b04f-6424-4b11-adb3-db1762f1a537.md"

$zhHandbackA = "0b47b63f-f910-4fa8-8c79-d4ab304f577c.e9b5a7c44d6d6801646440cd2948fdd45f83cb97.zh-cn.xlf"
$zhHandbackB = "c617b04f-6424-4b11-adb3-db1762f1a537.4f536b4b721dd6c2227b2eaae900273d212c7aa8.zh-cn.xlf"
$deHandbackA = "0b47b63f-f910-4fa8-8c79-d4ab304f577c.e9b5a7c44d6d6801646440cd2948fdd45f83cb97.de-de.xlf"
$deHandbackB = "c617b04f-6424-4b11-adb3-db1762f1a537.4f536b4b721dd6c2227b2eaae900273d212c7aa8.de-de.xlf"

$zhHandbackTime = "2016-09-05 18:54:09"
$deHandbackTime = "2016-09-05 18:54:18"

# Hyperlink font color used by the sheet's existing "HyperLink" style
# (blue FF6495ED, underlined) - OLE color is BGR, so 0x00ED9564.
$hyperlinkColor = 15570276

function Apply-HyperlinkLook($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: update the per-language status columns (E = zh-cn,
# F = de-de) for both rows to reflect the handback.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Columns widen to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column.
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Latest Target File -> hyperlink back to the (now in-sync) source file.
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, $null, $null, $srcMdA)
Apply-HyperlinkLook $wsZh.Range("I2")

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlB, $null, $null, $srcMdB)
Apply-HyperlinkLook $wsZh.Range("I3")

# Latest Handback File.
$wsZh.Range("J2").Value = $zhHandbackA
$wsZh.Range("J3").Value = $zhHandbackB

# Latest Handback DateTime.
$wsZh.Range("K2").Value = $zhHandbackTime
$wsZh.Range("K3").Value = $zhHandbackTime

# Widen the Status / Latest Target File / Latest Handback File columns.
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column.
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Latest Target File -> hyperlink back to the (now in-sync) source file.
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, $null, $null, $srcMdA)
Apply-HyperlinkLook $wsDe.Range("I2")

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlB, $null, $null, $srcMdB)
Apply-HyperlinkLook $wsDe.Range("I3")

# Latest Handback File.
$wsDe.Range("J2").Value = $deHandbackA
$wsDe.Range("J3").Value = $deHandbackB

# Latest Handback DateTime.
$wsDe.Range("K2").Value = $deHandbackTime
$wsDe.Range("K3").Value = $deHandbackTime

# Widen the Status / Latest Target File / Latest Handback File columns.
$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Output "Handback report generated."
